# Updated Time sheet - Sheik Fareeth
#
# Renames "Sheet1" -> "Day 7 (11-04-2022)", updates several timesheet rows
# (status / comments / hours) on that sheet, adjusts row heights to fit the
# new comment text, and refreshes the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rename the 5th sheet to reflect the new day ---------------------------
$ws.Name = "Day 7 (11-04-2022)"

# --- Row 8: no content change, only the wrapped comment needs more room ----
$ws.Rows.Item(8).RowHeight = 103.5

# --- Row 9: no content change, only the wrapped comment needs more room ----
$ws.Rows.Item(9).RowHeight = 99

# --- Row 10: comment text revised -------------------------------------------
$ws.Range("D10").Value = "1 Hr 30 mins : Brainstorming with team" + [char]10 + "1 Hr 30 mins : Prototype for TAC-(Create,View & Delete Pool -  6 slides)" + [char]10 + "1 Hr         : Redefined the overall prototype for TAC"
$ws.Rows.Item(10).RowHeight = 102.75

# --- Row 11: hours spent - project increased, taller comment ---------------
$ws.Range("F11").Value = 3.5
$ws.Rows.Item(11).RowHeight = 116.25

# --- Row 12: only row height shrinks ----------------------------------------
$ws.Rows.Item(12).RowHeight = 133.5

# --- Row 13: only row height shrinks ----------------------------------------
$ws.Rows.Item(13).RowHeight = 133.5

# --- Row 14: status changed to "Data Model" ---------------------------------
$ws.Range("C14").Value = "Data Model"

# --- Row 15: status + comment rewritten, hours increased -------------------
$ws.Range("C15").Value = "Data Model"
$ws.Range("D15").Value = "1 hr - Brain Stromming" + [char]10 + "30 Mins - Adding New slides to TAC ( Upcoming drives, notifications, scheduling drives - 5 slides )" + [char]10 + "2 hr - Building Data model In draw.io ( 7 Entities )" + [char]10 + "30 Mins - Re refining Interviewers scheduled drive cancellation"
$ws.Range("F15").Value = 4
$ws.Rows.Item(15).RowHeight = 144

# --- Row 16: only row height shrinks ----------------------------------------
$ws.Rows.Item(16).RowHeight = 111

# --- Row 17: status changed, hours spent - project decreased ---------------
$ws.Range("C17").Value = "Data Model"
$ws.Range("F17").Value = 3.5
$ws.Rows.Item(17).RowHeight = 142.5

# --- Refresh the active selection on the sheet ------------------------------
$ws.Activate()
$ws.Range("J10").Select()
